$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A85").Value = "GRT-USD"
